$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), copying the formatting from
# the existing header cell H1 (bold, centered, thin border - style index 1)
# so no duplicate style entry is created.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0 / IF data columns for rows 2-41
$iVals = 6, 6, 8, 6, 8, 8, 7, 9, 9, 8, 8, 7, 9, 7, 9, 8, 6, 6, 8, 7, 8, 9, 7, 8, 6, 7, 7, 8, 7, 8, 6, 5, 5, 8, 5, 6, 9, 8, 8, 5
$jVals = 6, 6, 9, 6, 8, 8, 7, 9, 9, 9, 8, 7, 9, 7, 9, 8, 6, 6, 8, 7, 8, 9, 7, 8, 6, 7, 7, 8, 7, 8, 6, 5, 6, 8, 5, 6, 9, 8, 8, 5

for ($i = 0; $i -lt $iVals.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $iVals[$i]
    $ws.Cells.Item($row, 10).Value = $jVals[$i]
}
